$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text for the MODEL_CONDITION column (currently in column E,
# before the column shift) to MODELCONDITION.
$ws.Range("E1").Value = "MODELCONDITION"

# Delete the (unlabeled) row-index column A entirely; this shifts
# columns B:F left to A:E, matching the target layout.
$ws.Columns.Item(1).Delete()
